# Optimize database schema migration
# Appends one new log row (row 48) to each of the four worksheets,
# mirroring the existing row layout (time, length, ID, actual length,
# checksum + their decimal counterparts).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW35-FE-LIFTER ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A48").Value = 45748.34298797454
$ws1.Range("A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B48").Value = "0x01,0x90"
$ws1.Range("C48").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Range("D48").Value = "0x01,0x76"
$ws1.Range("E48").Value = "0xd"
$ws1.Range("F48").Value = 400
$ws1.Range("G48").Value = 568631262647114.0 * [Math]::Pow(10, 9)
$ws1.Range("H48").Value = 374
$ws1.Range("I48").Value = 13

# --- Sheet 2: ROW35-MID-LIFTER ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A48").Value = 45748.19588700232
$ws2.Range("A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B48").Value = "0x01,0x90"
$ws2.Range("C48").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Range("D48").Value = "0x01,0x76"
$ws2.Range("E48").Value = "0xe"
$ws2.Range("F48").Value = 400
$ws2.Range("G48").Value = 568631262647114.0 * [Math]::Pow(10, 9)
$ws2.Range("H48").Value = 374
$ws2.Range("I48").Value = 14

# --- Sheet 3: ROW02-FE-LIFTER ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A48").Value = 45748.33583146991
$ws3.Range("A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B48").Value = "0x01,0x90"
$ws3.Range("C48").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Range("D48").Value = "0x01,0x76"
$ws3.Range("E48").Value = "0x3"
$ws3.Range("F48").Value = 400
$ws3.Range("G48").Value = 568631262647114.0 * [Math]::Pow(10, 9)
$ws3.Range("H48").Value = 374
$ws3.Range("I48").Value = 3

# --- Sheet 4: ROW02-MID-LIFTER ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A48").Value = 45748.3906115625
$ws4.Range("A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B48").Value = "0x01,0x90"
$ws4.Range("C48").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Range("D48").Value = "0x01,0x76"
$ws4.Range("E48").Value = "0x3"
$ws4.Range("F48").Value = 400
$ws4.Range("G48").Value = 985046333984776.0 * [Math]::Pow(10, 9)
$ws4.Range("H48").Value = 374
$ws4.Range("I48").Value = 3
